# Iselin WHOI 2016-11 Winklers — add the 111416 data sheet.
#
# Mirrors the structure of the existing "111116" sheet (same column
# layout / number formats) and appends it as the new last tab, which
# becomes the active one.

$wb = $excel.ActiveWorkbook

# The existing "111116" sheet (4th tab) is our formatting template.
$src = $wb.Worksheets.Item(4)

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "111416"

# Clone the cell formatting (wrap text / number format "0.000") from the
# template sheet's A1:H16 block onto the new sheet in one shot so the new
# sheet reuses the workbook's existing two cell styles exactly.
$src.Range("A1:H16").Copy()
$newSheet.Range("A1:H16").PasteSpecial(-4122)

# Seed the brand-new shared strings (station labels L29-L32) in the exact
# order they are first used (right-to-left across row 5-8's H column) so
# they land in the shared string table in the same order as the source
# workbook: L32, L31, L30, L29.
$newSheet.Range("H8").Value = "L32"
$newSheet.Range("H7").Value = "L31"
$newSheet.Range("H6").Value = "L30"
$newSheet.Range("H5").Value = "L29"

# Fill in the rest of the Winkler titration data table.
$newSheet.Range("A1").Value = "Phorcys"
$newSheet.Range("B1").Value = 23
$newSheet.Range("C1").Value = "T"
$newSheet.Range("D1").Value = 1
$newSheet.Range("E1").Value = 6.113
$newSheet.Range("F1").Value = "mL/L"
$newSheet.Range("G1").Value = "Oxygen"
$newSheet.Range("H1").Value = "T1"
$newSheet.Range("A2").Value = "Phorcys"
$newSheet.Range("B2").Value = 23
$newSheet.Range("C2").Value = "T"
$newSheet.Range("D2").Value = 2
$newSheet.Range("E2").Value = 6.172
$newSheet.Range("F2").Value = "mL/L"
$newSheet.Range("G2").Value = "Oxygen"
$newSheet.Range("H2").Value = "T2"
$newSheet.Range("A3").Value = "Phorcys"
$newSheet.Range("B3").Value = 23
$newSheet.Range("C3").Value = "T"
$newSheet.Range("D3").Value = 3
$newSheet.Range("E3").Value = 6.048
$newSheet.Range("F3").Value = "mL/L"
$newSheet.Range("G3").Value = "Oxygen"
$newSheet.Range("H3").Value = "T3"
$newSheet.Range("A4").Value = "Phorcys"
$newSheet.Range("B4").Value = 23
$newSheet.Range("C4").Value = "T"
$newSheet.Range("D4").Value = 4
$newSheet.Range("E4").Value = 6.324
$newSheet.Range("F4").Value = "mL/L"
$newSheet.Range("G4").Value = "Oxygen"
$newSheet.Range("H4").Value = "T4"
$newSheet.Range("A5").Value = "Phorcys"
$newSheet.Range("B5").Value = 24
$newSheet.Range("C5").Value = "L"
$newSheet.Range("D5").Value = 1
$newSheet.Range("E5").Value = 6.148
$newSheet.Range("F5").Value = "mL/L"
$newSheet.Range("G5").Value = "Oxygen"
$newSheet.Range("A6").Value = "Phorcys"
$newSheet.Range("B6").Value = 24
$newSheet.Range("C6").Value = "L"
$newSheet.Range("D6").Value = 2
$newSheet.Range("E6").Value = 6.164
$newSheet.Range("F6").Value = "mL/L"
$newSheet.Range("G6").Value = "Oxygen"
$newSheet.Range("A7").Value = "Phorcys"
$newSheet.Range("B7").Value = 24
$newSheet.Range("C7").Value = "L"
$newSheet.Range("D7").Value = 3
$newSheet.Range("E7").Value = 6.187
$newSheet.Range("F7").Value = "mL/L"
$newSheet.Range("G7").Value = "Oxygen"
$newSheet.Range("A8").Value = "Phorcys"
$newSheet.Range("B8").Value = 24
$newSheet.Range("C8").Value = "L"
$newSheet.Range("D8").Value = 4
$newSheet.Range("E8").Value = 6.172
$newSheet.Range("F8").Value = "mL/L"
$newSheet.Range("G8").Value = "Oxygen"
$newSheet.Range("A9").Value = "Phorcys"
$newSheet.Range("B9").Value = 25
$newSheet.Range("C9").Value = "F"
$newSheet.Range("D9").Value = 1
$newSheet.Range("E9").Value = 6.295
$newSheet.Range("F9").Value = "mL/L"
$newSheet.Range("G9").Value = "Oxygen"
$newSheet.Range("H9").Value = "F15"
$newSheet.Range("A10").Value = "Phorcys"
$newSheet.Range("B10").Value = 25
$newSheet.Range("C10").Value = "F"
$newSheet.Range("D10").Value = 2
$newSheet.Range("E10").Value = 6.35
$newSheet.Range("F10").Value = "mL/L"
$newSheet.Range("G10").Value = "Oxygen"
$newSheet.Range("H10").Value = "F16"
$newSheet.Range("A11").Value = "Phorcys"
$newSheet.Range("B11").Value = 25
$newSheet.Range("C11").Value = "F"
$newSheet.Range("D11").Value = 3
$newSheet.Range("E11").Value = 6.296
$newSheet.Range("F11").Value = "mL/L"
$newSheet.Range("G11").Value = "Oxygen"
$newSheet.Range("H11").Value = "F17"
$newSheet.Range("A12").Value = "Phorcys"
$newSheet.Range("B12").Value = 25
$newSheet.Range("C12").Value = "F"
$newSheet.Range("D12").Value = 4
$newSheet.Range("E12").Value = 6.303
$newSheet.Range("F12").Value = "mL/L"
$newSheet.Range("G12").Value = "Oxygen"
$newSheet.Range("H12").Value = "F18"
$newSheet.Range("A13").Value = "Phorcys"
$newSheet.Range("B13").Value = 26
$newSheet.Range("C13").Value = "F"
$newSheet.Range("D13").Value = 1
$newSheet.Range("E13").Value = 6.318
$newSheet.Range("F13").Value = "mL/L"
$newSheet.Range("G13").Value = "Oxygen"
$newSheet.Range("H13").Value = "F19"
$newSheet.Range("A14").Value = "Phorcys"
$newSheet.Range("B14").Value = 26
$newSheet.Range("C14").Value = "F"
$newSheet.Range("D14").Value = 2
$newSheet.Range("E14").Value = 6.532
$newSheet.Range("F14").Value = "mL/L"
$newSheet.Range("G14").Value = "Oxygen"
$newSheet.Range("H14").Value = "F20"
$newSheet.Range("A15").Value = "Phorcys"
$newSheet.Range("B15").Value = 26
$newSheet.Range("C15").Value = "F"
$newSheet.Range("D15").Value = 3
$newSheet.Range("E15").Value = 6.285
$newSheet.Range("F15").Value = "mL/L"
$newSheet.Range("G15").Value = "Oxygen"
$newSheet.Range("H15").Value = "F21"
$newSheet.Range("A16").Value = "Phorcys"
$newSheet.Range("B16").Value = 26
$newSheet.Range("C16").Value = "F"
$newSheet.Range("D16").Value = 4
$newSheet.Range("E16").Value = 6.318
$newSheet.Range("F16").Value = "mL/L"
$newSheet.Range("G16").Value = "Oxygen"
$newSheet.Range("H16").Value = "F22"

# Match the saved selection/active-cell state from the source file.
$newSheet.Range("I15").Select()
